{"js": "const body = context.document.body;\n\n// 1) \"megapolitan\" -> \"metapolitan\"\nconst hit1 = body.search(\"megapolitan\", { matchCase: true });\nhit1.load(\"items\");\nawait context.sync();\nif (hit1.items.length > 0) {\n  hit1.items[0].insertText(\"metapolitan\", \"Replace\");\n}\n\n// 2) \"investment banks)\" -> \"investment bankers)\"\nconst hit2 = body.search(\"investment banks)\", { matchCase: true });\nhit2.load(\"items\");\nawait context.sync();\nif (hit2.items.length > 0) {\n  hit2.items[0].insertText(\"investment bankers)\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Fix spelling: \"megapolitan\" -> \"metapolitan\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"megapolitan\"\n$find1.Replacement.Text = \"metapolitan\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, 2) | Out-Null\n\n# Fix spelling: \"banks\" -> \"bankers\" (in \"...top investment banks).\")\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"investment banks)\"\n$find2.Replacement.Text = \"investment bankers)\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2) | Out-Null\n"}
